# Add the new LeetCode entry "42. Trapping Rain Water" as row 15 of the
# summary table (category "2 Pointers"), matching the author's commit
# "did the trap rainwater problem".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: category (reuse existing "2 Pointers" shared string from row 12-14)
$ws.Range("A15").Value = "2 Pointers"

# Column B: question title
$ws.Range("B15").Value = "42. Trapping Rain Water"

# Column C: solution summary (long write-up, includes an embedded newline)
$ws.Range("C15").Value = "The idea is to use 2 ptrs, `"leftWall`" & `"rightWall`", use a while loop, inside call a func findRightWall(), this function finds the right wall give and arr & leftWall, the right wall is such that It is >= left wall, but incase there is no such wall in subarr leftWall to end, then return the biggest wall from subarr leftwall to end of arr.... next is a func getWaterBetween Walls() given an height arr, leftWall & rightWall, it returns the am of water that can be collected between them, use a while loop & initalize tmpPtr = leftWall+1, in while(tmpPtr<rightWall) loop do totalWater += height[tmpPtr] - minWallHeight where minWallHeight = min(height[leftWall],height[rightWall]), then tmpPter++`nFinally after calling both these funcs change the right wall to left wall ie leftWall = rightWall,,,, now in next iter we will find right wall again"

# Highlight the new question title cell with a red fill (new fill/style)
$ws.Range("B15").Interior.Color = 255

# Row height grows to fit the long wrapped text (8 lines @ 14.4pt)
$ws.Rows.Item(15).RowHeight = 115.2

# Update the active selection to reflect where the author finished editing
$ws.Range("E14").Select()
